# finish scene quest magnetic
# Append the new "magnetic" quest reward entry (42000014;1) to the two
# scene reward-list cells (F4 and F5), and move the active selection
# from F5 to F4 to match where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scene")

$ws.Range("F4").Value = "42000003;3|42000002;1|42000004;2|42000005;1|42000011;1|42000014;1"
$ws.Range("F5").Value = "42000006;2|42000007;1|42000008;2|42000003;3|42000004;2|42000013;1|42000014;1"

$ws.Range("F4").Select()
